$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: stash formats of source cells into a scratch area (row 30) before clearing ---
$ws.Range("C4").Copy($ws.Range("A30"))
$ws.Range("D4").Copy($ws.Range("B30"))
$ws.Range("E4").Copy($ws.Range("C30"))
$ws.Range("F4").Copy($ws.Range("D30"))
$ws.Range("C5").Copy($ws.Range("E30"))
$ws.Range("D5").Copy($ws.Range("F30"))
$ws.Range("E5").Copy($ws.Range("G30"))
$ws.Range("F5").Copy($ws.Range("H30"))
$ws.Range("B6").Copy($ws.Range("I30"))
$ws.Range("C6").Copy($ws.Range("J30"))
$ws.Range("D6").Copy($ws.Range("K30"))
$ws.Range("E6").Copy($ws.Range("L30"))
$ws.Range("F6").Copy($ws.Range("M30"))
$ws.Range("B7").Copy($ws.Range("N30"))
$ws.Range("C7").Copy($ws.Range("O30"))
$ws.Range("D7").Copy($ws.Range("P30"))
$ws.Range("E7").Copy($ws.Range("Q30"))
$ws.Range("F7").Copy($ws.Range("R30"))
$ws.Range("C8").Copy($ws.Range("S30"))
$ws.Range("D8").Copy($ws.Range("T30"))
$ws.Range("E8").Copy($ws.Range("U30"))
$ws.Range("F8").Copy($ws.Range("V30"))

# --- Step 2: clear the original used area only (keep scratch row 30 intact) ---
$ws.Range("A1:G9").Clear()

# --- Step 3: set column widths (COM ColumnWidth quantizes to MDW-7 pixel grid, closest achievable) ---
$ws.Columns.Item(1).ColumnWidth = 12.022135416666666
$ws.Columns.Item(2).ColumnWidth = 10.736979166666666
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 14.592447916666666
$ws.Columns.Item(5).ColumnWidth = 11.166666666666666

# --- Step 4: paste formats from scratch cells into final positions ---
$ws.Range("A30").Copy($ws.Range("B1"))
$ws.Range("B30").Copy($ws.Range("C1"))
$ws.Range("C30").Copy($ws.Range("D1"))
$ws.Range("D30").Copy($ws.Range("E1"))
$ws.Range("E30").Copy($ws.Range("B2"))
$ws.Range("F30").Copy($ws.Range("C2"))
$ws.Range("G30").Copy($ws.Range("D2"))
$ws.Range("H30").Copy($ws.Range("E2"))
$ws.Range("I30").Copy($ws.Range("A3"))
$ws.Range("J30").Copy($ws.Range("B3"))
$ws.Range("K30").Copy($ws.Range("C3"))
$ws.Range("L30").Copy($ws.Range("D3"))
$ws.Range("M30").Copy($ws.Range("E3"))
$ws.Range("N30").Copy($ws.Range("A4"))
$ws.Range("O30").Copy($ws.Range("B4"))
$ws.Range("P30").Copy($ws.Range("C4"))
$ws.Range("Q30").Copy($ws.Range("D4"))
$ws.Range("R30").Copy($ws.Range("E4"))
$ws.Range("S30").Copy($ws.Range("B5"))
$ws.Range("T30").Copy($ws.Range("C5"))
$ws.Range("U30").Copy($ws.Range("D5"))
$ws.Range("V30").Copy($ws.Range("E5"))

# --- Step 5: clear scratch area ---
$ws.Rows.Item(30).Clear()

# --- Step 6: set final values/text ---
$ws.Range("A1").Value = "Rownames"
$ws.Range("B1").Value = "Win_Low"
$ws.Range("C1").Value = "Win_MediumLow"
$ws.Range("D1").Value = "Win_MediumHigh"
$ws.Range("E1").Value = "Win_High"
$ws.Range("A2").Value = "Mac_Low"
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 17
$ws.Range("D2").Value = 17
$ws.Range("E2").Value = 5
$ws.Range("A3").Value = "Mac_MediumLow"
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = 22
$ws.Range("E3").Value = 27
$ws.Range("A4").Value = "Mac_MediumHigh"
$ws.Range("B4").Value = 16
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 14
$ws.Range("E4").Value = 20
$ws.Range("A5").Value = "Mac_High"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 17
$ws.Range("D5").Value = 19
$ws.Range("E5").Value = 11

# --- Step 7: selection to match target (B15) ---
$ws.Range("B15").Select()
